# SSD_CDS.xlsx update: rename series "DMS" -> "SSD", restate several historical
# monthly return figures in column B, and append 12 new monthly observations
# (Jan-2021 .. Dec-2021, rows 232-243).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text: "DMS" -> "SSD" (B1) ---
$ws.Range("B1").Value = "SSD"

# --- 2. Restate existing monthly return values in column B ---
# Map of row number -> new numeric value (as strings so we can safely parse
# scientific-notation literals with [double]::Parse / the [double] cast).
$changedB = @{
    12 = "0.14200343615821079"
    18 = "3.4453567888967918E-2"
    19 = "-5.5677953625615367E-2"
    20 = "3.1596488985554859E-2"
    21 = "3.1764135204536408E-2"
    38 = "4.4794022885418013E-2"
    41 = "7.3087728545139036E-2"
    42 = "-0.23353044697992609"
    45 = "-9.704731558160011E-2"
    46 = "-1.3011924959425641E-2"
    47 = "4.9657260937253873E-2"
    49 = "4.1737761762189997E-2"
    56 = "9.5602548519657216E-2"
    65 = "-4.3655561634834049E-2"
    66 = "-8.1992530803822822E-2"
    68 = "7.5112544194156786E-3"
    69 = "1.087221591966045E-2"
    70 = "-0.16850574394440901"
    71 = "2.7795741135660301E-2"
    80 = "3.7804483344180337E-2"
    82 = "-5.7452121343344528E-3"
    83 = "-2.1734019026089711E-2"
    84 = "-5.9913873768985673E-2"
    85 = "0.188208520361165"
    86 = "6.827542465373948E-2"
    87 = "-0.14382336696050169"
    88 = "6.9889262897853952E-2"
    89 = "6.5822285183019844E-2"
    90 = "3.7152224488639092E-2"
    91 = "-3.7334653062835732E-2"
    92 = "-2.1934775772478529E-2"
    93 = "5.358672442314652E-2"
    94 = "-7.9490592503852375E-2"
    95 = "3.0981108360163479E-2"
    96 = "3.5607213249702667E-2"
    98 = "4.6145787010186352E-2"
    120 = "-5.540297686022716E-2"
    175 = "2.4540523309063231E-2"
    202 = "-2.0500796221520859E-2"
    203 = "-8.5622753369916788E-3"
    204 = "-7.7762062386555961E-3"
    205 = "-8.2998242333865879E-2"
    206 = "-6.8039470079538811E-2"
    207 = "-1.6281413604064401E-2"
    208 = "3.5338147360052322E-2"
    209 = "3.2311913122180597E-2"
    210 = "5.8128320688908969E-2"
    211 = "2.2810098705751819E-2"
    212 = "-4.9860718927033212E-2"
    213 = "2.990523139750241E-2"
    222 = "2.3517257588701209E-2"
    223 = "0.11884213773735371"
    225 = "-8.0370107960096818E-2"
    226 = "8.8653886549407618E-2"
    227 = "-1.6248606744852609E-2"
    228 = "-4.3590935486754273E-2"
    229 = "1.1359668158807399E-2"
    230 = "1.8853798589539419E-2"
    231 = "7.3976386161182614E-2"
}
foreach ($r in $changedB.Keys) {
    $ws.Range("B" + $r).Value = [double]$changedB[$r]
}

# --- 3. Append new monthly rows 232-243 (date serial in col A, return in col B) ---
$newRows = @(
    @(232, "44225", "1.023640821955971E-2"),
    @(233, "44253", "9.0682419567015685E-2"),
    @(234, "44286", "-4.4689023254423171E-2"),
    @(235, "44316", "1.7652325248135151E-2"),
    @(236, "44347", "1.2142553017812771E-2"),
    @(237, "44377", "3.4352320730842059E-2"),
    @(238, "44407", "-4.721817998621116E-2"),
    @(239, "44439", "-1.6209740387514309E-2"),
    @(240, "44469", "3.7761984900447333E-2"),
    @(241, "44498", "2.9822591026699999E-2"),
    @(242, "44530", "-7.527781553960336E-2"),
    @(243, "44561", "3.1338824567149932E-2"),
)
foreach ($item in $newRows) {
    $r = $item[0]
    $dateSerial = [double]$item[1]
    $val = [double]$item[2]

    # Copy the date cell formatting (style) from the last existing data row so the
    # new dates keep the same "yyyy-mm-dd hh:mm:ss" display as the rest of column A.
    $ws.Range("A231").Copy() | Out-Null
    $ws.Range("A" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $r).Value = $dateSerial
    $ws.Range("B" + $r).Value = $val
}
$excel.CutCopyMode = 0

# --- 4. Give column B an explicit width now that it holds real data (col A already ---
# --- has a customWidth) ---
$ws.Columns.Item(2).ColumnWidth = [double]"11.619791666666666"

# --- 5. Dimension grows to A1:B243 automatically as cells are populated; update the ---
# --- active selection to match the saved view (J5) ---
$ws.Range("J5").Select() | Out-Null

Write-Output "SSD_CDS update applied"
